$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "force_error" column (F): clear the header label and all data values,
# leaving behind an empty, still-formatted F1 cell (as Excel does on Delete).
$ws.Range("F1:F25").ClearContents()

# Rename the "S2_Stabilizing" treatment label to "Drug" throughout column B.
$ws.Range("B8:B13").Value = "Drug"
$ws.Range("B20:B25").Value = "Drug"

# Leave the selection where the user finished editing.
$ws.Range("B25").Select()
